# Applies the "Updated cryptos list" data refresh to Sheet1 (columns B-E, rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to keep a numeric-looking Price (column D) entry as
# text -- exactly like typing `577.43 into the cell -- instead of auto-converting it
# to a number and dropping significant trailing/leading zeros.
$q = [char]39

$ws.Range("D2").Value = "61.626.85"
$ws.Range("E2").Value = "  +1.05%  "

$ws.Range("D3").Value = "3.390.61"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = $q + "577.43"
$ws.Range("E5").Value = "  +0.94%  "

$ws.Range("D6").Value = $q + "136.97"
$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "3.391.48"
$ws.Range("E8").Value = "  +0.61%  "

$ws.Range("E9").Value = "  -0.63%  "

$ws.Range("D10").Value = $q + "7.49"
$ws.Range("E10").Value = "  -1.45%  "

$ws.Range("D11").Value = $q + "0.126"
$ws.Range("E11").Value = "  +2.16%  "

$ws.Range("E12").Value = "  +0.63%  "

$ws.Range("D13").Value = "3.965.69"
$ws.Range("E13").Value = "  +0.26%  "

$ws.Range("E14").Value = "  +1.70%  "

$ws.Range("E15").Value = "  +1.77%  "

$ws.Range("D16").Value = "3.387.07"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").Value = $q + "25.67"
$ws.Range("E17").Value = "  +2.30%  "

$ws.Range("D18").Value = "61.733.70"
$ws.Range("E18").Value = "  +0.85%  "

$ws.Range("D19").Value = $q + "14.19"
$ws.Range("E19").Value = "  +1.26%  "

$ws.Range("E20").Value = "  +1.24%  "

$ws.Range("E21").Value = "  +0.34%  "

$ws.Range("D22").Value = $q + "378.15"
$ws.Range("E22").Value = "  +1.15%  "

$ws.Range("D23").Value = $q + "0.561"
$ws.Range("E23").Value = "  -1.20%  "

$ws.Range("D24").Value = "3.522.92"
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("E26").Value = "  +6.26%  "

$ws.Range("D27").Value = $q + "71.14"
$ws.Range("E27").Value = "  +0.73%  "

$ws.Range("D28").Value = $q + "1.75"
$ws.Range("E28").Value = "  +4.64%  "

$ws.Range("D29").Value = $q + "7.60"
$ws.Range("E29").Value = "  -0.99%  "

$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = $q + "0.160"
$ws.Range("E31").Value = "  +3.58%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = $q + "8.19"
$ws.Range("E32").Value = "  +1.01%  "

$ws.Range("E33").Value = "  +0.78%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Value = $q + "23.37"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").Value = $q + "5.35"
$ws.Range("E36").Value = "  -3.62%  "

$ws.Range("E37").Value = "  -0.61%  "

$ws.Range("D38").Value = $q + "6.85"
$ws.Range("E38").Value = "  -0.90%  "

$ws.Range("D39").Value = $q + "164.89"
$ws.Range("E39").Value = "  +1.23%  "

$ws.Range("D40").Value = $q + "0.0786"
$ws.Range("E40").Value = "  -0.34%  "

$ws.Range("D41").Value = $q + "0.782"
$ws.Range("E41").Value = "  +2.96%  "

$ws.Range("B42").Value = "ONDO"
$ws.Range("C42").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D42").Value = $q + "1.24"
$ws.Range("E42").Value = "  +2.56%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = $q + "1.00"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("E44").Value = "  +7.67%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = $q + "24.95"
$ws.Range("E45").Value = "  +7.72%  "

$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = $q + "4.41"
$ws.Range("E46").Value = "  -0.08%  "

$ws.Range("D47").Value = $q + "41.34"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("E48").Value = "  -1.70%  "

$ws.Range("D49").Value = $q + "22.74"
$ws.Range("E49").Value = "  -1.00%  "

$ws.Range("D50").Value = "2.336.18"
$ws.Range("E50").Value = "  +5.89%  "

$ws.Range("E51").Value = "  +1.29%  "
